$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was updated from
# 2023-10-22 (45221) to 2023-10-25 (45224) for every data row (2-99).
$ws.Range("C2:C99").Value = 45224
